$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$sub3 = [char]0x2083   # unicode subscript three, used in ShibaInu price

Set-TextValue "D2" "42.092.71"
Set-TextValue "E2" "  -1.93%  "
Set-TextValue "D3" "2.245.65"
Set-TextValue "E3" "  -1.88%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "247.55"
Set-TextValue "E5" "  -1.93%  "
Set-TextValue "D6" "0.625"
Set-TextValue "E6" "  -2.66%  "
Set-TextValue "D7" "77.09"
Set-TextValue "E7" "  +4.37%  "
Set-TextValue "E8" "  +0.06%  "
Set-TextValue "D9" "0.627"
Set-TextValue "E9" "  -3.28%  "
Set-TextValue "D10" "41.85"
Set-TextValue "E10" "  +7.51%  "
Set-TextValue "D11" "0.0957"
Set-TextValue "E11" "  -2.39%  "
Set-TextValue "E12" "  -2.14%  "
Set-TextValue "E13" "  -3.18%  "
Set-TextValue "D14" "2.580.24"
Set-TextValue "E14" "  -1.95%  "
Set-TextValue "E15" "  -2.81%  "
Set-TextValue "D16" "0.862"
Set-TextValue "E16" "  -1.52%  "
Set-TextValue "D17" "2.220.97"
Set-TextValue "E17" "  -3.14%  "
Set-TextValue "D18" "42.003.02"
Set-TextValue "E18" "  -1.93%  "
Set-TextValue "D19" ("0.0{0}0984" -f $sub3)
Set-TextValue "E19" "  -2.28%  "
Set-TextValue "E20" "  -2.68%  "
Set-TextValue "D21" "72.00"
Set-TextValue "E21" "  -1.13%  "
Set-TextValue "D22" "2.33"
Set-TextValue "E22" "  +4.81%  "
Set-TextValue "D23" "231.85"
Set-TextValue "E23" "  -1.39%  "
Set-TextValue "B24" "Dai"
Set-TextValue "C24" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D24" "1.00"
Set-TextValue "E24" "  +0.04%  "
Set-TextValue "B25" "Cosmos"
Set-TextValue "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D25" "11.40"
Set-TextValue "E25" "  -1.52%  "
Set-TextValue "D26" "3.66"
Set-TextValue "E26" "  -6.13%  "
Set-TextValue "D27" "2.30"
Set-TextValue "E27" "  -5.38%  "
Set-TextValue "D28" "7.29"
Set-TextValue "E28" "  +12.46%  "
Set-TextValue "E29" "  +0.82%  "
Set-TextValue "D30" "169.43"
Set-TextValue "E30" "  +1.47%  "
Set-TextValue "D31" "20.58"
Set-TextValue "E31" "  -2.33%  "
Set-TextValue "D32" "33.67"
Set-TextValue "E32" "  +8.16%  "
Set-TextValue "D33" "0.0831"
Set-TextValue "E33" "  +0.23%  "
Set-TextValue "E34" "  -4.98%  "
Set-TextValue "E35" "  -0.79%  "
Set-TextValue "D36" "4.57"
Set-TextValue "E36" "  -1.18%  "
Set-TextValue "E37" "  +2.93%  "
Set-TextValue "B38" "Celestia"
Set-TextValue "C38" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D38" "14.36"
Set-TextValue "E38" "  -0.25%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0302"
Set-TextValue "E39" "  -2.45%  "
Set-TextValue "E40" "  +0.20%  "
Set-TextValue "E41" "  -6.55%  "
Set-TextValue "D42" "113.28"
Set-TextValue "E42" "  +11.70%  "
Set-TextValue "E43" "  -6.10%  "
Set-TextValue "D44" "61.29"
Set-TextValue "E44" "  -1.03%  "
Set-TextValue "E45" "  -4.79%  "
Set-TextValue "E46" "  -2.74%  "
Set-TextValue "D47" "0.997"
Set-TextValue "E47" "  -0.38%  "
Set-TextValue "E48" "  -2.95%  "
Set-TextValue "E49" "  -1.14%  "
Set-TextValue "D50" "4.23"
Set-TextValue "E50" "  -13.27%  "
Set-TextValue "D51" "2.29"
Set-TextValue "E51" "  -0.70%  "
